# Updates "Price" (D) and "Volume(1h)" (E) columns for the latest cryptos
# snapshot. Price values are stored as plain text in the workbook (not
# numbers), so numeric-looking prices are entered with a leading apostrophe
# to force text entry, then the cell style is reset to "Normal" so no
# stray number-format / quote-prefix style lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.596.40"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.76%  "
$ws.Range("D3").Value = "'2.509.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.30%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'575.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.76%  "
$ws.Range("D6").Value = "'166.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.76%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "'2.508.25"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.32%  "
$ws.Range("E10").Value = "  -7.14%  "
$ws.Range("D12").Value = "'0.342"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.08%  "
$ws.Range("D13").Value = "'4.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.47%  "
$ws.Range("D14").Value = "'2.969.42"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.27%  "
$ws.Range("D15").Value = "'69.519.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("E16").Value = "  -6.36%  "
$ws.Range("D17").Value = "'24.88"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.55%  "
$ws.Range("D18").Value = "'2.510.73"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.55%  "
$ws.Range("D19").Value = "'11.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.95%  "
$ws.Range("D20").Value = "'7.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").Value = "'350.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.94%  "
$ws.Range("D22").Value = "'3.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.25%  "
$ws.Range("E23").Value = "  -5.27%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "'68.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.67%  "
$ws.Range("D26").Value = "'4.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.46%  "
$ws.Range("D27").Value = "'8.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -7.54%  "
$ws.Range("D28").Value = "'2.639.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'0.0₃0902"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.42%  "
$ws.Range("D31").Value = "'7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("D32").Value = "'476.58"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.37%  "
$ws.Range("E33").Value = "  -0.73%  "
$ws.Range("E34").Value = "  -3.54%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").Value = "'154.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.56%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("D39").Value = "'18.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.31%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").Value = "'0.319"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.09%  "
$ws.Range("E43").Value = "  -7.85%  "
$ws.Range("E44").Value = "  -13.31%  "
$ws.Range("D45").Value = "'2.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.95%  "
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").Value = "'144.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.29%  "
$ws.Range("D48").Value = "'0.530"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.09%  "
$ws.Range("E49").Value = "  -3.98%  "
$ws.Range("E50").Value = "  -5.51%  "
$ws.Range("E51").Value = "  -2.70%  "
